$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for the two additional series (columns H and I)
$ws.Range("H1").Value = "MA.TDA_AL_multiTier_TDAamortAS_OYLM_TDA_LowRate2_base"
$ws.Range("I1").Value = "MA.TDA_AL_multiTier_TDAamortAS_OYLM_TDA_LowRate2_lowG"

# New data values for rows 2-4
$ws.Range("H2").Value = 0.2722982533414307
$ws.Range("I2").Value = 0.2722982533414307

$ws.Range("H3").Value = 0.2507626990404917
$ws.Range("I3").Value = 0.21304163431459305

$ws.Range("H4").Value = 0.2200181755456755
$ws.Range("I4").Value = 0.1363696517206612
